$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: "002" -> "001" (must stay text, not become the number 1, so build it via a
# text formula and paste-special the computed value back in as a literal).
$ws.Range("J2").Formula = '="001"'
$ws.Range("J2").Copy()
$ws.Range("J2").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("N2").Value = "2017-12-31 00:00:00"

$ws.Range("O2").Value = 743203940.23
$ws.Range("P2").Value = 269245209.19
$ws.Range("Q2").Value = 65907719.84
$ws.Range("S2").Value = 190194449.76
$ws.Range("U2").Value = 90685822.15000001
$ws.Range("W2").Value = 530967737.62
$ws.Range("X2").Value = 150292522.79
$ws.Range("Z2").Value = 5646952.63
$ws.Range("AB2").Value = 212236202.61
$ws.Range("AF2").Value = 72.0340731336
$ws.Range("AG2").Value = 71.4430735466
